$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Deletions (processed top-to-bottom by original row number so the
#     later, larger row indices are still valid when we get to them) ---

# Row 9: 004384167 / DOUGLAS / 29698.94 - removed entirely
$ws.Rows.Item(9).Delete()

# After the delete above, the two stale "MIRELLA" rows (originally rows
# 100 and 101: 003553997/200.03 and 001651617/200.02) have shifted up to
# rows 99 and 100. Remove them (same index twice, since deleting row 99
# shifts the old row 100 into row 99).
$ws.Rows.Item(99).Delete()
$ws.Rows.Item(99).Delete()

# --- Insertions ---
# New rows re-appear higher up the sheet with updated balances. Insert the
# deeper one (before 004290978/LARISSA, now at row 15) first so it doesn't
# disturb the row number of the shallower insertion point.
$ws.Rows.Item(15).Insert()
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = "001651617"
$ws.Cells.Item(15, 2).Value = "MIRELLA"
$ws.Cells.Item(15, 3).Value = 2000

# Insert before 001759765/NATAL (row 13 after the earlier delete).
$ws.Rows.Item(13).Insert()
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "003553997"
$ws.Cells.Item(13, 2).Value = "MIRELLA"
$ws.Cells.Item(13, 3).Value = 5000
